$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, pushing existing rows 29-44 down to 30-45.
$ws.Rows.Item(29).Insert()

# Fill in the new record's data (matches the other rows' shared values for most columns).
$ws.Cells.Item(29, 1).Value = 5
$ws.Cells.Item(29, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(29, 3).Value = "Maule"
$ws.Cells.Item(29, 4).Value = 44784
$ws.Cells.Item(29, 5).Value = 7
$ws.Cells.Item(29, 6).Value = 100112040
$ws.Cells.Item(29, 7).Value = "Cilantro"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 15000
$ws.Cells.Item(29, 12).Value = 15000
$ws.Cells.Item(29, 13).Value = 15000
$ws.Cells.Item(29, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(29, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(29, 16).Value = 417
$ws.Cells.Item(29, 17).Value = 36
$ws.Cells.Item(29, 18).Value = "Hortaliza"
